$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete duplicate observation rows. Deleting bottom-up keeps the
# remaining row numbers stable while each Delete() executes.
$ws.Rows("21").Delete()
$ws.Rows("19").Delete()
$ws.Rows("18").Delete()
$ws.Rows("15").Delete()
$ws.Rows("13").Delete()

# The rows that remain (now 12-16) lose their Report_Volume /
# Observation_Date values and get re-flagged as "Neutral" ENSO entries.
for ($r = 12; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).ClearContents()
    $ws.Cells.Item($r, 2).ClearContents()
    $c = $ws.Cells.Item($r, 3)
    $c.Value = "Neutral"
    $c.Font.Size = 12
}

$ws.Range("D12").Select() | Out-Null
